# Update the Inputs sheet's "Capital Cost" value (Inputs!C2) from
# 400,000,000 to 100,000,000. All downstream formulas (Simluation_Tool_Example
# and Outputs sheets) recalculate automatically from this single input change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inputs")
$ws.Activate()

$ws.Range("C2").Select()
$ws.Range("C2").Value = 100000000

$ws.Range("C3").Select()
